{"js": "// Update the division-fact answers in the single results table.\n// Each data row (0, 4, 8, 12, 16 -- the other rows are spacer rows) has\n// 5 cells; we replace the text of each cell in place while preserving\n// the existing run/paragraph formatting (font, size, alignment) by\n// targeting the paragraph's Range and doing a \"replace\" insert instead\n// of rewriting the whole cell body.\n\nconst replacements = [\n  [0, 0, \"25\u00f78=3, 1\", \"29\u00f77=4, 1\"],\n  [0, 1, \"65\u00f79=7, 2\", \"90\u00f76=15, 0\"],\n  [0, 2, \"34\u00f76=5, 4\", \"60\u00f75=12, 0\"],\n  [0, 3, \"91\u00f73=30, 1\", \"84\u00f77=12, 0\"],\n  [0, 4, \"62\u00f77=8, 6\", \"84\u00f75=16, 4\"],\n  [4, 0, \"46\u00f73=15, 1\", \"70\u00f76=11, 4\"],\n  [4, 1, \"68\u00f75=13, 3\", \"33\u00f72=16, 1\"],\n  [4, 2, \"81\u00f75=16, 1\", \"36\u00f77=5, 1\"],\n  [4, 3, \"67\u00f77=9, 4\", \"59\u00f76=9, 5\"],\n  [4, 4, \"59\u00f78=7, 3\", \"83\u00f72=41, 1\"],\n  [8, 0, \"53\u00f78=6, 5\", \"62\u00f78=7, 6\"],\n  [8, 1, \"95\u00f76=15, 5\", \"22\u00f74=5, 2\"],\n  [8, 2, \"35\u00f79=3, 8\", \"84\u00f78=10, 4\"],\n  [8, 3, \"85\u00f74=21, 1\", \"76\u00f79=8, 4\"],\n  [8, 4, \"59\u00f77=8, 3\", \"84\u00f74=21, 0\"],\n  [12, 0, \"85\u00f74=21, 1\", \"46\u00f77=6, 4\"],\n  [12, 1, \"27\u00f73=9, 0\", \"76\u00f78=9, 4\"],\n  [12, 2, \"61\u00f78=7, 5\", \"12\u00f74=3, 0\"],\n  [12, 3, \"74\u00f72=37, 0\", \"61\u00f73=20, 1\"],\n  [12, 4, \"33\u00f77=4, 5\", \"64\u00f77=9, 1\"],\n  [16, 0, \"33\u00f72=16, 1\", \"36\u00f76=6, 0\"],\n  [16, 1, \"98\u00f77=14, 0\", \"10\u00f75=2, 0\"],\n  [16, 2, \"41\u00f72=20, 1\", \"98\u00f72=49, 0\"],\n  [16, 3, \"62\u00f75=12, 2\", \"70\u00f73=23, 1\"],\n  [16, 4, \"51\u00f75=10, 1\", \"40\u00f74=10, 0\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\n// Grab every cell we need to touch up front.\nconst cells = replacements.map(([row, col]) => table.getCell(row, col));\nfor (const cell of cells) {\n  cell.body.paragraphs.load(\"items\");\n}\nawait context.sync();\n\nconst ranges = cells.map((cell) => cell.body.paragraphs.items[0].getRange());\nfor (const range of ranges) {\n  range.load(\"text\");\n}\nawait context.sync();\n\nfor (let i = 0; i < replacements.length; i++) {\n  const [, , oldText, newText] = replacements[i];\n  const range = ranges[i];\n  if (range.text !== oldText) {\n    throw new Error(\n      `Unexpected cell text \"${range.text}\" (expected \"${oldText}\") at index ${i}`\n    );\n  }\n  range.insertText(newText, Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Update the division-fact answers in the single results table.\n# The table has 20 rows x 5 columns; only rows 1, 5, 9, 13, 17 (1-based)\n# hold data -- the rows in between are blank spacer rows. We overwrite\n# each data cell's Range.Text in place, which keeps the existing\n# paragraph/run formatting (font, size, alignment) untouched.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$replacements = @(\n    @{ Row = 1;  Col = 1; Old = \"25\u00f78=3, 1\";   New = \"29\u00f77=4, 1\" },\n    @{ Row = 1;  Col = 2; Old = \"65\u00f79=7, 2\";   New = \"90\u00f76=15, 0\" },\n    @{ Row = 1;  Col = 3; Old = \"34\u00f76=5, 4\";   New = \"60\u00f75=12, 0\" },\n    @{ Row = 1;  Col = 4; Old = \"91\u00f73=30, 1\";  New = \"84\u00f77=12, 0\" },\n    @{ Row = 1;  Col = 5; Old = \"62\u00f77=8, 6\";   New = \"84\u00f75=16, 4\" },\n\n    @{ Row = 5;  Col = 1; Old = \"46\u00f73=15, 1\";  New = \"70\u00f76=11, 4\" },\n    @{ Row = 5;  Col = 2; Old = \"68\u00f75=13, 3\";  New = \"33\u00f72=16, 1\" },\n    @{ Row = 5;  Col = 3; Old = \"81\u00f75=16, 1\";  New = \"36\u00f77=5, 1\" },\n    @{ Row = 5;  Col = 4; Old = \"67\u00f77=9, 4\";   New = \"59\u00f76=9, 5\" },\n    @{ Row = 5;  Col = 5; Old = \"59\u00f78=7, 3\";   New = \"83\u00f72=41, 1\" },\n\n    @{ Row = 9;  Col = 1; Old = \"53\u00f78=6, 5\";   New = \"62\u00f78=7, 6\" },\n    @{ Row = 9;  Col = 2; Old = \"95\u00f76=15, 5\";  New = \"22\u00f74=5, 2\" },\n    @{ Row = 9;  Col = 3; Old = \"35\u00f79=3, 8\";   New = \"84\u00f78=10, 4\" },\n    @{ Row = 9;  Col = 4; Old = \"85\u00f74=21, 1\";  New = \"76\u00f79=8, 4\" },\n    @{ Row = 9;  Col = 5; Old = \"59\u00f77=8, 3\";   New = \"84\u00f74=21, 0\" },\n\n    @{ Row = 13; Col = 1; Old = \"85\u00f74=21, 1\";  New = \"46\u00f77=6, 4\" },\n    @{ Row = 13; Col = 2; Old = \"27\u00f73=9, 0\";   New = \"76\u00f78=9, 4\" },\n    @{ Row = 13; Col = 3; Old = \"61\u00f78=7, 5\";   New = \"12\u00f74=3, 0\" },\n    @{ Row = 13; Col = 4; Old = \"74\u00f72=37, 0\";  New = \"61\u00f73=20, 1\" },\n    @{ Row = 13; Col = 5; Old = \"33\u00f77=4, 5\";   New = \"64\u00f77=9, 1\" },\n\n    @{ Row = 17; Col = 1; Old = \"33\u00f72=16, 1\";  New = \"36\u00f76=6, 0\" },\n    @{ Row = 17; Col = 2; Old = \"98\u00f77=14, 0\";  New = \"10\u00f75=2, 0\" },\n    @{ Row = 17; Col = 3; Old = \"41\u00f72=20, 1\";  New = \"98\u00f72=49, 0\" },\n    @{ Row = 17; Col = 4; Old = \"62\u00f75=12, 2\";  New = \"70\u00f73=23, 1\" },\n    @{ Row = 17; Col = 5; Old = \"51\u00f75=10, 1\";  New = \"40\u00f74=10, 0\" }\n)\n\nforeach ($r in $replacements) {\n    $cell = $t.Cell($r.Row, $r.Col)\n    $range = $cell.Range\n    # Cell range text includes the trailing end-of-cell marker(s); strip\n    # them off before comparing so we can sanity-check the current value.\n    $current = $range.Text.TrimEnd([char]13, [char]7)\n    if ($current -ne $r.Old) {\n        throw \"Unexpected cell text '$current' (expected '$($r.Old)') at row $($r.Row) col $($r.Col)\"\n    }\n    $range.Text = $r.New\n}\n"}
